# Update column C (Förändrad) for rows 2-20 from 45174 to 45175
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 20; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45174) {
        $cell.Value = 45175
    }
}
